# Re-applies the refreshed cryptos price/volume snapshot (GitHub Actions run).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells in the
# sheet (t="inlineStr"), even when their content looks numeric (e.g. "0.9997").
# A leading apostrophe forces Excel to store the new value as literal text too,
# instead of silently re-typing it as a Number (which would e.g. truncate
# "8.520" down to "8.52").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    if ($text -match '^[+-]?[0-9.]+$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue 'D2' '26.884.71'
Set-TextValue 'E2' '  -0.67%  '
Set-TextValue 'D3' '1.859.83'
Set-TextValue 'E3' '  -0.13%  '
Set-TextValue 'D4' '0.9994'
Set-TextValue 'D5' '304.56'
Set-TextValue 'E5' '  -0.55%  '
Set-TextValue 'D6' '0.9994'
Set-TextValue 'E6' '  -0.01%  '
Set-TextValue 'D7' '0.5037'
Set-TextValue 'E7' '  -1.61%  '
Set-TextValue 'E8' '  -2.60%  '
Set-TextValue 'D9' '0.07169'
Set-TextValue 'E9' '  +0.73%  '
Set-TextValue 'D10' '0.8934'
Set-TextValue 'E10' '  +0.31%  '
Set-TextValue 'E11' '  +0.48%  '
Set-TextValue 'D12' '1.869.62'
Set-TextValue 'E12' '  +0.47%  '
Set-TextValue 'D13' '0.07486'
Set-TextValue 'E13' '  -0.76%  '
Set-TextValue 'D14' '94.15'
Set-TextValue 'E14' '  +5.92%  '
Set-TextValue 'D15' '5.228'
Set-TextValue 'E15' '  -1.38%  '
Set-TextValue 'D16' '0.9997'
Set-TextValue 'E16' '  -0.04%  '
Set-TextValue 'D17' '0.000008497'
Set-TextValue 'E17' '  +1.68%  '
Set-TextValue 'E18' '  +0.68%  '
Set-TextValue 'D19' '0.9992'
Set-TextValue 'E19' '  -0.05%  '
Set-TextValue 'D20' '26.938.99'
Set-TextValue 'E20' '  -0.60%  '
Set-TextValue 'D21' '5.021'
Set-TextValue 'E21' '  -0.80%  '
Set-TextValue 'D22' '2.107.88'
Set-TextValue 'E22' '  +0.84%  '
Set-TextValue 'D24' '6.413'
Set-TextValue 'E24' '  -0.87%  '
Set-TextValue 'D25' '147.89'
Set-TextValue 'E25' '  -0.87%  '
Set-TextValue 'D26' '1.778'
Set-TextValue 'E26' '  -3.41%  '
Set-TextValue 'E27' '  -0.67%  '
Set-TextValue 'E28' '  -0.13%  '
Set-TextValue 'D29' '112.99'
Set-TextValue 'E29' '  +0.08%  '
Set-TextValue 'D30' '4.692'
Set-TextValue 'E30' '  +0.08%  '
Set-TextValue 'D31' '4.669'
Set-TextValue 'E31' '  +0.17%  '
Set-TextValue 'D32' '0.09221'
Set-TextValue 'E32' '  +2.04%  '
Set-TextValue 'D33' '0.05145'
Set-TextValue 'E33' '  +0.58%  '
Set-TextValue 'D34' '0.7470'
Set-TextValue 'D35' '2.966'
Set-TextValue 'E35' '  -2.68%  '
Set-TextValue 'D36' '1.151'
Set-TextValue 'E36' '  -0.58%  '
Set-TextValue 'E37' '  +6.43%  '
Set-TextValue 'D38' '2.578'
Set-TextValue 'E38' '  +2.92%  '
Set-TextValue 'D39' '0.02006'
Set-TextValue 'E39' '  -2.10%  '
Set-TextValue 'D40' '0.5553'
Set-TextValue 'E40' '  +4.03%  '
Set-TextValue 'E41' '  -0.05%  '
Set-TextValue 'D42' '6.554'
Set-TextValue 'E42' '  -0.56%  '
Set-TextValue 'D43' '117.35'
Set-TextValue 'E43' '  +1.30%  '
Set-TextValue 'D44' '8.520'
Set-TextValue 'E44' '  +2.44%  '
Set-TextValue 'E45' '  -0.09%  '
Set-TextValue 'D46' '0.4678'
Set-TextValue 'E46' '  +1.13%  '
Set-TextValue 'D47' '0.9991'
Set-TextValue 'E47' '  -0.02%  '
Set-TextValue 'D48' '10.01'
Set-TextValue 'E48' '  -0.50%  '
Set-TextValue 'D49' '1.563'
Set-TextValue 'E49' '  -0.14%  '
Set-TextValue 'E50' '  -0.43%  '
Set-TextValue 'D51' '62.98'
Set-TextValue 'E51' '  -1.79%  '
